# Updates market/leve profit figures (columns H-N) across several sheets,
# reflecting refreshed price data pulled by the scheduled runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 47
$ws.Range("H47").Value = 0
$ws.Range("I47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("K47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("M47").Value = ""
$ws.Range("N47").Value = ""

# Row 48
$ws.Range("H48").Value = 2000
$ws.Range("I48").Value = 0
$ws.Range("K48").Value = 0
$ws.Range("M48").Value = ""

# Row 56
$ws.Range("H56").Value = 2000
$ws.Range("I56").Value = 0
$ws.Range("K56").Value = 0
$ws.Range("M56").Value = ""

# Row 135
$ws.Range("H135").Value = 1010.9231
$ws.Range("I135").Value = 890.7
$ws.Range("K135").Value = 8016.3
$ws.Range("M135").Value = -5481.3

# Row 137
$ws.Range("H137").Value = 10096.667
$ws.Range("I137").Value = 3309.32
$ws.Range("K137").Value = 9927.960000000001
$ws.Range("M137").Value = -7377.960000000001

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 4316.7085
$ws.Range("I32").Value = 2786.9023
$ws.Range("J32").Value = 13277
$ws.Range("K32").Value = 2786.9023
$ws.Range("L32").Value = 13277
$ws.Range("M32").Value = -2499.9023
$ws.Range("N32").Value = -13851

# Row 61
$ws.Range("H61").Value = 1138837.1
$ws.Range("I61").Value = 3635.7727
$ws.Range("K61").Value = 3635.7727
$ws.Range("M61").Value = -3423.7727

# Row 74
$ws.Range("H74").Value = 23916.3
$ws.Range("I74").Value = 4953.933
$ws.Range("J74").Value = 80803.39999999999
$ws.Range("K74").Value = 4953.933
$ws.Range("L74").Value = 80803.39999999999
$ws.Range("M74").Value = -4079.933
$ws.Range("N74").Value = -82551.39999999999

# Row 77
$ws.Range("H77").Value = 23916.3
$ws.Range("I77").Value = 4953.933
$ws.Range("J77").Value = 80803.39999999999
$ws.Range("K77").Value = 24769.665
$ws.Range("L77").Value = 404017
$ws.Range("M77").Value = -20401.665
$ws.Range("N77").Value = -412753

# Row 109
$ws.Range("H109").Value = 87500
$ws.Range("J109").Value = 87500
$ws.Range("L109").Value = 87500
$ws.Range("N109").Value = -90274

# Row 121
$ws.Range("H121").Value = 218000
$ws.Range("J121").Value = 218000
$ws.Range("L121").Value = 218000
$ws.Range("N121").Value = -221494

# Row 132
$ws.Range("H132").Value = 14720895
$ws.Range("I132").Value = 5992.2
$ws.Range("K132").Value = 17976.6
$ws.Range("M132").Value = -15446.6

# Row 136
$ws.Range("H136").Value = 1138837.1
$ws.Range("I136").Value = 3635.7727
$ws.Range("K136").Value = 10907.3181
$ws.Range("M136").Value = -8357.3181

$ws = $wb.Worksheets.Item("BSM")
# Row 81
$ws.Range("H81").Value = 43583.625
$ws.Range("J81").Value = 43583.625
$ws.Range("L81").Value = 43583.625
$ws.Range("N81").Value = -45705.625

# Row 84
$ws.Range("H84").Value = 43583.625
$ws.Range("J84").Value = 43583.625
$ws.Range("L84").Value = 130750.875
$ws.Range("N84").Value = -141358.875

# Row 103
$ws.Range("H103").Value = 15998.5
$ws.Range("J103").Value = 15998.5
$ws.Range("L103").Value = 15998.5
$ws.Range("N103").Value = -18342.5

# Row 134
$ws.Range("H134").Value = 11149.617
$ws.Range("I134").Value = 5951.317
$ws.Range("K134").Value = 17853.951
$ws.Range("M134").Value = -15318.951

$ws = $wb.Worksheets.Item("CRP")
# Row 5
$ws.Range("H5").Value = 1706
$ws.Range("I5").Value = 307.8889
$ws.Range("J5").Value = 7997.5
$ws.Range("K5").Value = 307.8889
$ws.Range("L5").Value = 7997.5
$ws.Range("M5").Value = -195.8889
$ws.Range("N5").Value = -8221.5

# Row 31
$ws.Range("H31").Value = 45694.516
$ws.Range("I31").Value = 98188.62
$ws.Range("K31").Value = 98188.62
$ws.Range("M31").Value = -97893.62

# Row 34
$ws.Range("H34").Value = 45694.516
$ws.Range("I34").Value = 98188.62
$ws.Range("K34").Value = 98188.62
$ws.Range("M34").Value = -97986.62

# Row 105
$ws.Range("H105").Value = 12418
$ws.Range("I105").Value = 13345.25
$ws.Range("K105").Value = 13345.25
$ws.Range("M105").Value = -11598.25

# Row 134
$ws.Range("H134").Value = 28576988
$ws.Range("J134").Value = 90923020
$ws.Range("L134").Value = 272769060
$ws.Range("N134").Value = -272774130

$ws = $wb.Worksheets.Item("CUL")
# Row 23
$ws.Range("H23").Value = 178.8
$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("M23").Value = ""

# Row 49
$ws.Range("H49").Value = 596.5
$ws.Range("I49").Value = 193
$ws.Range("K49").Value = 579
$ws.Range("M49").Value = -423

# Row 61
$ws.Range("H61").Value = 139.41667
$ws.Range("I61").Value = 171.11111
$ws.Range("J61").Value = 44.333332
$ws.Range("K61").Value = 513.3333299999999
$ws.Range("L61").Value = 132.999996
$ws.Range("M61").Value = -298.3333299999999
$ws.Range("N61").Value = -562.999996

$ws = $wb.Worksheets.Item("GSM")
# Row 132
$ws.Range("H132").Value = 1320097.2
$ws.Range("I132").Value = 9709.799999999999
$ws.Range("K132").Value = 29129.4
$ws.Range("M132").Value = -26599.4

$ws = $wb.Worksheets.Item("LTW")
# Row 68
$ws.Range("H68").Value = 53333.332
$ws.Range("I68").Value = 53333.332
$ws.Range("K68").Value = 53333.332
$ws.Range("M68").Value = -52584.332

# Row 71
$ws.Range("H71").Value = 53333.332
$ws.Range("I71").Value = 53333.332
$ws.Range("K71").Value = 266666.66
$ws.Range("M71").Value = -262922.66

$ws = $wb.Worksheets.Item("WVR")
# Row 29
$ws.Range("H29").Value = 50000
$ws.Range("J29").Value = 50000
$ws.Range("L29").Value = 50000
$ws.Range("N29").Value = -50580

# Row 62
$ws.Range("H62").Value = 56667.332
$ws.Range("I62").Value = 50000
$ws.Range("J62").Value = 58000.8
$ws.Range("K62").Value = 50000
$ws.Range("L62").Value = 58000.8
$ws.Range("M62").Value = -49376
$ws.Range("N62").Value = -59248.8

# Row 65
$ws.Range("H65").Value = 56667.332
$ws.Range("I65").Value = 50000
$ws.Range("J65").Value = 58000.8
$ws.Range("K65").Value = 250000
$ws.Range("L65").Value = 290004
$ws.Range("M65").Value = -246880
$ws.Range("N65").Value = -296244

# Row 122
$ws.Range("H122").Value = 5050.3335
$ws.Range("I122").Value = 3067.5
$ws.Range("J122").Value = 7033.1665
$ws.Range("K122").Value = 9202.5
$ws.Range("L122").Value = 21099.4995
$ws.Range("M122").Value = -6752.5
$ws.Range("N122").Value = -25999.4995

# Row 136
$ws.Range("H136").Value = 614388.9
$ws.Range("J136").Value = 1488446.2
$ws.Range("L136").Value = 4465338.6
$ws.Range("N136").Value = -4470438.6
